# Add a "Month" column (L) with the invoice month, formatted with a
# month/year date-style number format, and update the selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new column.
$ws.Range("L1").Value = "Month"

# Data cell: the month name, styled with a date number format (mmm-yy).
$ws.Range("L2").Value = "February"
$ws.Range("L2").NumberFormat = "mmm-yy"

# Widen the new column to fit its content.
$ws.Columns(12).ColumnWidth = 16.28515625

# Move the active selection, as recorded after the edit.
$ws.Range("I6").Select()
